# Update code for caclular Luong ca nhan
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn 1 bác sĩ")

# Row 2: service group changed from "Ngực" to "Tiểu phẫu"
$ws.Range("G2").Value = "Tiểu phẫu"

# Insert a new data row (row 3) for the new invoice, pushing the previous
# "Tổng" (total) row down to row 4 along with its existing (blank) cells.
$ws.Rows.Item(3).Insert()

# New row 3 data (new invoice record)
$ws.Range("A3").Value = "HD-LUXURY"
$ws.Range("B3").Value = 554
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "07-15-2024"
$ws.Range("D3").Value = "SÓC TRĂNG"
$ws.Range("E3").Value = "lê thị bích"
$ws.Range("F3").Value = "Cá nhân"
$ws.Range("G3").Value = "Đại phẫu"
$ws.Range("H3").Value = "cắt sẹo "
$ws.Range("I3").Value = "Lâm Thị Mỹ Hằng"
$ws.Range("J3").Value = 7000000
$ws.Range("K3").Value = "Lê Đình Hậu"
$ws.Range("L3").Value = 6000000
$ws.Range("M3").Value = 13000000
$ws.Range("N3").Value = 13000000
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 13000000
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = "Bác Sĩ Thảo"
$ws.Range("T3").Value = "Trần Khánh Hiệp"
$ws.Range("X3").Value = 0.1
$ws.Range("Y3").Value = 0.04
$ws.Range("Z3").Value = 1060000
$ws.Range("AA3").Value = 240000

# Row 4 ("Tổng" summary row) - update only the totals that changed; the
# remaining (already blank) cells were carried down automatically by the
# row insert above.
$ws.Range("B4").Value = 2
$ws.Range("J4").Value = 32000000
$ws.Range("L4").Value = 6000000
$ws.Range("M4").Value = 38000000
$ws.Range("N4").Value = 38000000
$ws.Range("P4").Value = 38000000
$ws.Range("X4").Value = 0.1
$ws.Range("Y4").Value = 0.04
$ws.Range("Z4").Value = 1060000
$ws.Range("AA4").Value = 240000
